# update target editor tooltip
# ---------------------------------------------------------------------------
# 1) "Data" sheet: three new rows of sample log data (row 2..4), column F
#    holds a blank/newline "editor tooltip" cell rendered in a monospace
#    (Consolas 10pt) font - this is the actual "editor tooltip" styling
#    referenced by the commit message.
# 2) "Username" sheet: the two user names that appear in the new Data rows.
# 3) "DailyTarget" sheet: the stray empty A1 placeholder cell is cleared so
#    only the "TARGET" header in B1 remains.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Data sheet ------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$rows = @(
    @{ Row = 2; Num = 1; Date = "2025-04-21"; Code = "LU18"; User = "gg" },
    @{ Row = 3; Num = 2; Date = "2025-04-21"; Code = "LU18"; User = "aa" },
    @{ Row = 4; Num = 3; Date = "2025-04-21"; Code = "LU18"; User = "aa" }
)

foreach ($r in $rows) {
    $data.Cells.Item($r.Row, 1).Value = $r.Num             # A: row number

    # Column B holds a literal yyyy-mm-dd label (not a real Excel date), so
    # mark it as Text to stop the autoconverter turning it into a date
    # serial, write the string, then drop the temporary number-format
    # override now that the text is safely stored (done row-by-row, before
    # anything else on the row is touched, so it cannot bleed into later
    # edits on the same row).
    $data.Cells.Item($r.Row, 2).NumberFormat = "@"
    $data.Cells.Item($r.Row, 2).Value = $r.Date
    $data.Cells.Item($r.Row, 2).ClearFormats()

    $data.Cells.Item($r.Row, 4).Value = $r.Code            # D: code
    $data.Cells.Item($r.Row, 5).Value = $r.User            # E: user
    $data.Cells.Item($r.Row, 6).Value = "`n"               # F: tooltip cell
}

# Column F is the "editor tooltip" cell - rendered with a small monospace
# font so multi-line notes line up.
$data.Range("F2:F4").Font.Name = "Consolas"
$data.Range("F2:F4").Font.Size = 10

# --- Username sheet ----------------------------------------------------
$username = $wb.Worksheets.Item("Username")
$username.Range("A1").Value = "gg"
$username.Range("A2").Value = "aa"

# --- DailyTarget sheet ---------------------------------------------------
$dailyTarget = $wb.Worksheets.Item("DailyTarget")
$dailyTarget.Range("A1").ClearContents()
